$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.Value = "'308.74"
$c.Style = 'Normal'

$c = $ws.Range('E2')
$c.Value = "'-0.65%"
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.Value = "'39.47"
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.Value = "'1.52%"
$c.Style = 'Normal'

$c = $ws.Range('D4')
$c.Value = "'5.134"
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.Value = "'0.29%"
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.Value = "'0.08131"
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.Value = "'-0.54%"
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.Value = "'-3.25%"
$c.Style = 'Normal'

$c = $ws.Range('B7')
$c.Value = "'GateToken"
$c.Style = 'Normal'

$c = $ws.Range('C7')
$c.Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c.Style = 'Normal'

$c = $ws.Range('D7')
$c.Value = "'4.234"
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.Value = "'0.93%"
$c.Style = 'Normal'

$c = $ws.Range('B8')
$c.Value = "'KuCoinToken"
$c.Style = 'Normal'

$c = $ws.Range('C8')
$c.Value = "'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$c.Style = 'Normal'

$c = $ws.Range('D8')
$c.Value = "'8.129"
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.Value = "'2.49%"
$c.Style = 'Normal'

$c = $ws.Range('D9')
$c.Value = "'0.9265"
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.Value = "'-0.74%"
$c.Style = 'Normal'

$c = $ws.Range('D10')
$c.Value = "'0.1414"
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.Value = "'0.40%"
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.Value = "'-1.76%"
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.Value = "'0.09071"
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.Value = "'-1.14%"
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.Value = "'0.03493"
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.Value = "'1.09%"
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.Value = "'0.09798"
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.Value = "'-0.44%"
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.Value = "'0.001391"
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.Value = "'-1.98%"
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.Value = "'0.005839"
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.Value = "'-0.85%"
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.Value = "'3.915"
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.Value = "'9.59%"
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.Value = "'0.3454"
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.Value = "'0.12%"
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.Value = "'0.1313"
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.Value = "'-0.01%"
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.Value = "'4.731"
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.Value = "'-2.15%"
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.Value = "'-1.69%"
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.Value = "'0.04388"
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.Value = "'-1.71%"
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.Value = "'0.001232"
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.Value = "'-0.34%"
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.Value = "'0.004866"
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.Value = "'16.59%"
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.Value = "'0.0001301"
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.Value = "'-0.03%"
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.Value = "'0.0004004"
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.Value = "'-9.96%"
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.Value = "'0.02070"
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.Value = "'-3.33%"
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.Value = "'0.05066"
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.Value = "'-2.43%"
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.Value = "'0.007436"
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.Value = "'-0.54%"
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.Value = "'0.009760"
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.Value = "'-2.06%"
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.Value = "'0.1366"
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.Value = "'-0.23%"
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.Value = "'0.002132"
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.Value = "'-0.03%"
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.Value = "'0.009563"
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.Value = "'-2.01%"
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.Value = "'0.00006424"
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.Value = "'1.40%"
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.Value = "'0.00000000751"
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.Value = "'-0.03%"
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.Value = "'-18.83%"
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.Value = "'0.00002102"
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.Value = "'-0.03%"
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.Value = "'0.0002002"
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.Value = "'-0.03%"
$c.Style = 'Normal'
